$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original values for the columns that change (D, H, J, K, L, M, P)
# for rows 2, 3, 5, 6, 7 before overwriting anything, since the update is a
# cyclic re-shuffle of rows: new2=old7, new3=old2, new7=old3 (3-cycle) and
# new5=old6, new6=old5 (swap).

function Get-RowData($row) {
    return @{
        D = $ws.Range("D$row").Value2
        H = $ws.Range("H$row").Value2
        J = $ws.Range("J$row").Value2
        K = $ws.Range("K$row").Value2
        L = $ws.Range("L$row").Value2
        M = $ws.Range("M$row").Value2
        P = $ws.Range("P$row").Value2
    }
}

function Set-RowData($row, $data) {
    $ws.Range("D$row").Value = $data.D
    $ws.Range("H$row").Value = $data.H
    $ws.Range("J$row").Value = $data.J
    $ws.Range("K$row").Value = $data.K
    $ws.Range("L$row").Value = $data.L
    $ws.Range("M$row").Value = $data.M
    $ws.Range("P$row").Value = $data.P
}

$row2 = Get-RowData 2
$row3 = Get-RowData 3
$row5 = Get-RowData 5
$row6 = Get-RowData 6
$row7 = Get-RowData 7

# 3-cycle among rows 2, 3, 7: new2 = old7, new3 = old2, new7 = old3
Set-RowData 2 $row7
Set-RowData 3 $row2
Set-RowData 7 $row3

# swap rows 5 and 6
Set-RowData 5 $row6
Set-RowData 6 $row5
